# Applies the "automatic update" edit to the BJUV worksheet:
#  - every row's "Förändrad" (column C) date advances by one day (46065 -> 46066)
#  - rows 7..25 get their Beteckning/Datum/Area values reshuffled to reflect the
#    refreshed data export (including the Markägare note moving from row 9 to row 25)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C ("Förändrad") bumps by one day for every data row (2..25) ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = 46066
}

# --- Row-by-row content for the rows that changed (Beteckning / Datum / Area) ---
$rows = @{
    7  = @{ A = "A 59877-2025"; B = 45993;               G = 1 }
    9  = @{ A = "A 57394-2024"; B = 45629.6907175926;    G = 0.5 }
    10 = @{ A = "A 46993-2025"; B = 45929.54670138889;   G = 2.8 }
    11 = @{ A = "A 46998-2025"; B = 45929.54851851852;   G = 0.9 }
    12 = @{ A = "A 59471-2024"; B = 45638;               G = 1.8 }
    13 = @{ A = "A 20239-2025"; B = 45772;               G = 1.9 }
    14 = @{ A = "A 18118-2022"; B = 44684;               G = 4.2 }
    15 = @{ A = "A 57391-2024"; B = 45629.68717592592;   G = 1.8 }
    16 = @{ A = "A 54557-2023"; B = 45233.6346875;       G = 3.8 }
    18 = @{ A = "A 53218-2023"; B = 45229;               G = 5.4 }
    19 = @{ A = "A 20054-2024"; B = 45434;               G = 7.3 }
    20 = @{ A = "A 26708-2023"; B = 45093;               G = 4.1 }
    21 = @{ A = "A 63548-2025"; B = 46013;               G = 0.9 }
    22 = @{ A = "A 6679-2026";  B = 46056.60961805555;   G = 2.4 }
    23 = @{ A = "A 6684-2026";  B = 46056.61989583333;   G = 8.199999999999999 }
    24 = @{ A = "A 14674-2024"; B = 45397;               G = 4.8 }
    25 = @{ A = "A 4574-2025";  B = 45687;               G = 6.2 }
}

foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 7).Value = $row.G
}

# --- Markägare (column F) note moves from row 9 to row 25 ---
$ws.Cells.Item(9, 6).Value = ""
$ws.Cells.Item(25, 6).Value = "Kyrkan"
